$wb = $excel.ActiveWorkbook
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    Write-Output "$i : $($ws.Name)"
}
